$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.103.40"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.810.10"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.323"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").Value = "1.807.75"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "11.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "35.068.44"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +19.60%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  -6.06%  "
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "92.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.680"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "1.312.30"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  +4.92%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "1.989.32"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0654"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.74%  "
